$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 290.84616
$ws.Cells.Item(9, 9).Value = 298.41666
$ws.Cells.Item(9, 11).Value = 298.41666
$ws.Cells.Item(9, 13).Value = -129.41666
$ws.Cells.Item(11, 8).Value = 34370.043
$ws.Cells.Item(11, 9).Value = 34370.043
$ws.Cells.Item(11, 11).Value = 34370.043
$ws.Cells.Item(11, 13).Value = -34230.043
$ws.Cells.Item(19, 8).Value = 966.4666999999999
$ws.Cells.Item(19, 10).Value = 1059.7
$ws.Cells.Item(19, 12).Value = 1059.7
$ws.Cells.Item(19, 14).Value = -1409.7
$ws.Cells.Item(28, 8).Value = 1058.9474
$ws.Cells.Item(28, 9).Value = 1095.4615
$ws.Cells.Item(28, 10).Value = 979.8333
$ws.Cells.Item(28, 11).Value = 1095.4615
$ws.Cells.Item(28, 12).Value = 979.8333
$ws.Cells.Item(28, 13).Value = -610.4614999999999
$ws.Cells.Item(28, 14).Value = -1949.8333
$ws.Cells.Item(40, 8).Value = 1295.5625
$ws.Cells.Item(40, 9).Value = 1133.2222
$ws.Cells.Item(40, 10).Value = 1504.2858
$ws.Cells.Item(40, 11).Value = 1133.2222
$ws.Cells.Item(40, 12).Value = 1504.2858
$ws.Cells.Item(40, 13).Value = -958.2221999999999
$ws.Cells.Item(40, 14).Value = -1854.2858
$ws.Cells.Item(43, 8).Value = 3208.4546
$ws.Cells.Item(43, 10).Value = 3484.8572
$ws.Cells.Item(43, 12).Value = 3484.8572
$ws.Cells.Item(43, 14).Value = -3622.8572
$ws.Cells.Item(62, 8).Value = 2127.5
$ws.Cells.Item(62, 9).Value = 2127.5
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 2127.5
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 13).Value = -1503.5
$ws.Cells.Item(62, 14).Value = $null
$ws.Cells.Item(64, 8).Value = 6999.7144
$ws.Cells.Item(64, 10).Value = 8332.666999999999
$ws.Cells.Item(64, 12).Value = 8332.666999999999
$ws.Cells.Item(64, 14).Value = -8828.666999999999
$ws.Cells.Item(65, 8).Value = 2127.5
$ws.Cells.Item(65, 9).Value = 2127.5
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 10637.5
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = -7517.5
$ws.Cells.Item(65, 14).Value = $null
$ws.Cells.Item(67, 8).Value = 6999.7144
$ws.Cells.Item(67, 10).Value = 8332.666999999999
$ws.Cells.Item(67, 12).Value = 8332.666999999999
$ws.Cells.Item(67, 14).Value = -10048.667
$ws.Cells.Item(94, 8).Value = 25740.75
$ws.Cells.Item(94, 9).Value = 18748.625
$ws.Cells.Item(94, 11).Value = 18748.625
$ws.Cells.Item(94, 13).Value = -18297.625
$ws.Cells.Item(112, 8).Value = 3454.8948
$ws.Cells.Item(112, 10).Value = 4288.6665
$ws.Cells.Item(112, 12).Value = 12865.9995
$ws.Cells.Item(112, 14).Value = -15081.9995
$ws.Cells.Item(113, 8).Value = 5199.75
$ws.Cells.Item(113, 9).Value = 3599.6667
$ws.Cells.Item(113, 11).Value = 3599.6667
$ws.Cells.Item(113, 13).Value = -345.6667000000002
$ws.Cells.Item(116, 8).Value = 27380.5
$ws.Cells.Item(116, 10).Value = 31658.25
$ws.Cells.Item(116, 12).Value = 31658.25
$ws.Cells.Item(116, 14).Value = -38542.25
$ws.Cells.Item(135, 8).Value = 1585.275
$ws.Cells.Item(135, 9).Value = 1430.0588
$ws.Cells.Item(135, 10).Value = 2464.8333
$ws.Cells.Item(135, 11).Value = 12870.5292
$ws.Cells.Item(135, 12).Value = 22183.4997
$ws.Cells.Item(135, 13).Value = -10335.5292
$ws.Cells.Item(135, 14).Value = -27253.4997
$ws.Cells.Item(138, 8).Value = 7368.7144
$ws.Cells.Item(138, 9).Value = 7510.96
$ws.Cells.Item(138, 11).Value = 22532.88
$ws.Cells.Item(138, 13).Value = -17392.88
$ws.Cells.Item(63, 8).Value = 45000
$ws.Cells.Item(63, 10).Value = 45000
$ws.Cells.Item(63, 12).Value = 45000
$ws.Cells.Item(63, 14).Value = -46372
$ws.Cells.Item(66, 8).Value = 45000
$ws.Cells.Item(66, 10).Value = 45000
$ws.Cells.Item(66, 12).Value = 135000
$ws.Cells.Item(66, 14).Value = -141864

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(22, 8).Value = 2604
$ws.Cells.Item(22, 9).Value = 1472
$ws.Cells.Item(22, 10).Value = 6000
$ws.Cells.Item(22, 11).Value = 1472
$ws.Cells.Item(22, 12).Value = 6000
$ws.Cells.Item(22, 13).Value = -1173
$ws.Cells.Item(22, 14).Value = -6598
$ws.Cells.Item(25, 8).Value = 5262.9
$ws.Cells.Item(25, 10).Value = 7276.4287
$ws.Cells.Item(25, 12).Value = 7276.4287
$ws.Cells.Item(25, 14).Value = -8080.4287
$ws.Cells.Item(32, 8).Value = 732896.0600000001
$ws.Cells.Item(32, 9).Value = 806465.0600000001
$ws.Cells.Item(32, 10).Value = 7716
$ws.Cells.Item(32, 11).Value = 806465.0600000001
$ws.Cells.Item(32, 12).Value = 7716
$ws.Cells.Item(32, 13).Value = -806178.0600000001
$ws.Cells.Item(32, 14).Value = -8290
$ws.Cells.Item(61, 8).Value = 2780544.2
$ws.Cells.Item(61, 9).Value = 2654.303
$ws.Cells.Item(61, 11).Value = 2654.303
$ws.Cells.Item(61, 13).Value = -2442.303
$ws.Cells.Item(102, 8).Value = 3957.4285
$ws.Cells.Item(102, 9).Value = 3144.0625
$ws.Cells.Item(102, 10).Value = 6560.2
$ws.Cells.Item(102, 11).Value = 3144.0625
$ws.Cells.Item(102, 12).Value = 6560.2
$ws.Cells.Item(102, 13).Value = -1522.0625
$ws.Cells.Item(102, 14).Value = -9804.200000000001
$ws.Cells.Item(110, 8).Value = 1239.8572
$ws.Cells.Item(110, 9).Value = 1247.6666
$ws.Cells.Item(110, 10).Value = 1193
$ws.Cells.Item(110, 11).Value = 1247.6666
$ws.Cells.Item(110, 12).Value = 1193
$ws.Cells.Item(110, 13).Value = 797.3334
$ws.Cells.Item(110, 14).Value = -5283
$ws.Cells.Item(132, 8).Value = 864454
$ws.Cells.Item(132, 9).Value = 928160.0600000001
$ws.Cells.Item(132, 11).Value = 2784480.18
$ws.Cells.Item(132, 13).Value = -2781950.18
$ws.Cells.Item(136, 8).Value = 2780544.2
$ws.Cells.Item(136, 9).Value = 2654.303
$ws.Cells.Item(136, 11).Value = 7962.909
$ws.Cells.Item(136, 13).Value = -5412.909

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 58777.05
$ws.Cells.Item(20, 9).Value = 84421.234
$ws.Cells.Item(20, 11).Value = 84421.234
$ws.Cells.Item(20, 13).Value = -84174.234
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 14).Value = $null
$ws.Cells.Item(86, 8).Value = 2476.1
$ws.Cells.Item(86, 9).Value = 2451
$ws.Cells.Item(86, 11).Value = 2451
$ws.Cells.Item(86, 13).Value = -1328
$ws.Cells.Item(89, 8).Value = 2476.1
$ws.Cells.Item(89, 9).Value = 2451
$ws.Cells.Item(89, 11).Value = 12255
$ws.Cells.Item(89, 13).Value = -6639
$ws.Cells.Item(134, 8).Value = 9269417
$ws.Cells.Item(134, 9).Value = 4408.2144
$ws.Cells.Item(134, 11).Value = 13224.6432
$ws.Cells.Item(134, 13).Value = -10689.6432

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(13, 8).Value = 225
$ws.Cells.Item(13, 10).Value = 225
$ws.Cells.Item(13, 12).Value = 225
$ws.Cells.Item(13, 14).Value = -503
$ws.Cells.Item(31, 8).Value = 773712.25
$ws.Cells.Item(31, 9).Value = 1159457.4
$ws.Cells.Item(31, 11).Value = 1159457.4
$ws.Cells.Item(31, 13).Value = -1159162.4
$ws.Cells.Item(34, 8).Value = 773712.25
$ws.Cells.Item(34, 9).Value = 1159457.4
$ws.Cells.Item(34, 11).Value = 1159457.4
$ws.Cells.Item(34, 13).Value = -1159255.4
$ws.Cells.Item(62, 8).Value = 5186.4287
$ws.Cells.Item(62, 9).Value = 5498
$ws.Cells.Item(62, 10).Value = 5134.5
$ws.Cells.Item(62, 11).Value = 5498
$ws.Cells.Item(62, 12).Value = 5134.5
$ws.Cells.Item(62, 13).Value = -4874
$ws.Cells.Item(62, 14).Value = -6382.5
$ws.Cells.Item(65, 8).Value = 5186.4287
$ws.Cells.Item(65, 9).Value = 5498
$ws.Cells.Item(65, 10).Value = 5134.5
$ws.Cells.Item(65, 11).Value = 27490
$ws.Cells.Item(65, 12).Value = 25672.5
$ws.Cells.Item(65, 13).Value = -24370
$ws.Cells.Item(65, 14).Value = -31912.5
$ws.Cells.Item(99, 8).Value = 18850.77
$ws.Cells.Item(99, 10).Value = 3061.375
$ws.Cells.Item(99, 12).Value = 3061.375
$ws.Cells.Item(99, 14).Value = -6057.375
$ws.Cells.Item(126, 8).Value = 18850.77
$ws.Cells.Item(126, 10).Value = 3061.375
$ws.Cells.Item(126, 12).Value = 9184.125
$ws.Cells.Item(126, 14).Value = -14124.125
$ws.Cells.Item(132, 8).Value = 1145.6086
$ws.Cells.Item(132, 9).Value = 1017.45
$ws.Cells.Item(132, 10).Value = 2000
$ws.Cells.Item(132, 11).Value = 3052.35
$ws.Cells.Item(132, 12).Value = 6000
$ws.Cells.Item(132, 13).Value = -522.3500000000004
$ws.Cells.Item(132, 14).Value = -11060
$ws.Cells.Item(134, 8).Value = 2410.1892
$ws.Cells.Item(134, 9).Value = 2217.1482
$ws.Cells.Item(134, 10).Value = 2931.4
$ws.Cells.Item(134, 11).Value = 6651.444600000001
$ws.Cells.Item(134, 12).Value = 8794.200000000001
$ws.Cells.Item(134, 13).Value = -4116.444600000001
$ws.Cells.Item(134, 14).Value = -13864.2

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(26, 8).Value = 319.77777
$ws.Cells.Item(26, 9).Value = 207.9
$ws.Cells.Item(26, 11).Value = 623.7
$ws.Cells.Item(26, 13).Value = -335.7
$ws.Cells.Item(132, 8).Value = 8883.117
$ws.Cells.Item(132, 9).Value = 3188.375
$ws.Cells.Item(132, 11).Value = 28695.375
$ws.Cells.Item(132, 13).Value = -26165.375
$ws.Cells.Item(136, 8).Value = 9276.846
$ws.Cells.Item(136, 10).Value = 11124.25
$ws.Cells.Item(136, 12).Value = 33372.75
$ws.Cells.Item(136, 14).Value = -43572.75
$ws.Cells.Item(140, 8).Value = 1431.3572
$ws.Cells.Item(140, 9).Value = 1431.3572
$ws.Cells.Item(140, 11).Value = 4294.071599999999
$ws.Cells.Item(140, 13).Value = 885.9284000000007
$ws.Cells.Item(141, 8).Value = 9583
$ws.Cells.Item(141, 9).Value = 8499.6
$ws.Cells.Item(141, 11).Value = 25498.8
$ws.Cells.Item(141, 13).Value = -20318.8

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 3236.5386
$ws.Cells.Item(113, 9).Value = 2915.9092
$ws.Cells.Item(113, 11).Value = 2915.9092
$ws.Cells.Item(113, 13).Value = -745.9092000000001
$ws.Cells.Item(122, 8).Value = 7078.269
$ws.Cells.Item(122, 9).Value = 8340.888999999999
$ws.Cells.Item(122, 10).Value = 4237.375
$ws.Cells.Item(122, 11).Value = 25022.667
$ws.Cells.Item(122, 12).Value = 12712.125
$ws.Cells.Item(122, 13).Value = -22572.667
$ws.Cells.Item(122, 14).Value = -17612.125
$ws.Cells.Item(123, 8).Value = 37360.855
$ws.Cells.Item(123, 10).Value = 37360.855
$ws.Cells.Item(123, 12).Value = 37360.855
$ws.Cells.Item(123, 14).Value = -42260.855

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 26250.875
$ws.Cells.Item(7, 9).Value = 28572.428
$ws.Cells.Item(7, 11).Value = 28572.428
$ws.Cells.Item(7, 13).Value = -28460.428
$ws.Cells.Item(122, 8).Value = 4711.1113
$ws.Cells.Item(122, 9).Value = 3771.4285
$ws.Cells.Item(122, 11).Value = 11314.2855
$ws.Cells.Item(122, 13).Value = -8864.2855
$ws.Cells.Item(126, 8).Value = 26250.875
$ws.Cells.Item(126, 9).Value = 28572.428
$ws.Cells.Item(126, 11).Value = 85717.284
$ws.Cells.Item(126, 13).Value = -83247.284
$ws.Cells.Item(132, 8).Value = 7410406
$ws.Cells.Item(132, 9).Value = 9526808
$ws.Cells.Item(132, 10).Value = 3000
$ws.Cells.Item(132, 11).Value = 28580424
$ws.Cells.Item(132, 12).Value = 9000
$ws.Cells.Item(132, 13).Value = -28577894
$ws.Cells.Item(132, 14).Value = -14060
$ws.Cells.Item(136, 8).Value = 62505348
$ws.Cells.Item(136, 9).Value = 41673800
$ws.Cells.Item(136, 10).Value = 125000000
$ws.Cells.Item(136, 11).Value = 125021400
$ws.Cells.Item(136, 12).Value = 375000000
$ws.Cells.Item(136, 13).Value = -125018850
$ws.Cells.Item(136, 14).Value = -375005100

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(45, 8).Value = 15751.5
$ws.Cells.Item(45, 9).Value = 17187
$ws.Cells.Item(45, 11).Value = 17187
$ws.Cells.Item(45, 13).Value = -16696
$ws.Cells.Item(113, 8).Value = 541.3333
$ws.Cells.Item(113, 9).Value = 433.26315
$ws.Cells.Item(113, 11).Value = 1299.78945
$ws.Cells.Item(113, 13).Value = 870.21055
$ws.Cells.Item(132, 8).Value = 41670390
$ws.Cells.Item(132, 9).Value = 166666670
$ws.Cells.Item(132, 10).Value = 4966.6665
$ws.Cells.Item(132, 11).Value = 500000010
$ws.Cells.Item(132, 12).Value = 14899.9995
$ws.Cells.Item(132, 13).Value = -499997480
$ws.Cells.Item(132, 14).Value = -19959.9995
$ws.Cells.Item(136, 8).Value = 14350284
$ws.Cells.Item(136, 9).Value = 6213692.5
$ws.Cells.Item(136, 10).Value = 33335664
$ws.Cells.Item(136, 11).Value = 18641077.5
$ws.Cells.Item(136, 12).Value = 100006992
$ws.Cells.Item(136, 13).Value = -18638527.5
$ws.Cells.Item(136, 14).Value = -100012092
